$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.289.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.803.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.54%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("E6").Value = "  -4.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.802.81"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.57%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.30%  "
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000258"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.436.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.802.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.338.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "515.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.38%  "
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000139"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +25.09%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  +3.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.341"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("E40").Value = "  +3.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "44.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "423.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.068.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.76%  "
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "135.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("E51").Value = "  -0.43%  "
